$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Feria Lagunitas de Puerto Montt - Albahaca) needs
# to be inserted before the current row 95, pushing the existing rows 95-98
# down to 96-99.
$ws.Rows.Item(95).Insert()

# Populate the newly inserted row 95 with the new weekly record.
$ws.Range("A95").Value = 4
$ws.Range("B95").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C95").Value = "Los Lagos"
$ws.Range("D95").Value = 44610
$ws.Range("E95").Value = 10
$ws.Range("F95").Value = 100112052
$ws.Range("G95").Value = "Albahaca"
$ws.Range("H95").Value = "Sin especificar"
$ws.Range("I95").Value = "Primera"
$ws.Range("J95").Value = 150
$ws.Range("K95").Value = 6000
$ws.Range("L95").Value = 6000
$ws.Range("M95").Value = 6000
$ws.Range("N95").Value = "$/docena de matas"
$ws.Range("O95").Value = "Región Metropolitana"
$ws.Range("P95").Value = 1000
$ws.Range("Q95").Value = 6
$ws.Range("R95").Value = "Hortaliza"
